$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: Test_Pattern -> Plant_Type
$ws.Range("D1").Value = "Plant_Type"

# A new "Asphalt" sample row is introduced, shifting the existing
# Palm/Broad Leaf/Piney/Big Leaf/Coral labels down by one row, while the
# old "Piny Green" row is dropped. Update the Plant_Type (column D) labels
# accordingly.
$ws.Range("D2").Value = "Asphalt"
$ws.Range("D3").Value = "Palm"
$ws.Range("D4").Value = "Broad Leaf"
$ws.Range("D5").Value = "Piney"
$ws.Range("D6").Value = "Big Leaf"
$ws.Range("D7").Value = "Coral"

# Updated calibration measurements (Red_Endmember, NIR_Endmember, NDVI_Endmember)
$ws.Range("A2").Value = 0.05644740811392837
$ws.Range("B2").Value = 0.07865813722574855
$ws.Range("C2").Value = 0.164395392180083

$ws.Range("A3").Value = 0.1878047200867781
$ws.Range("B3").Value = 0.3882576717741545
$ws.Range("C3").Value = 0.3479709047484005

$ws.Range("A4").Value = 0.09359334001739344
$ws.Range("B4").Value = 1.108918777358624
$ws.Range("C4").Value = 0.8443369698068043

$ws.Range("A5").Value = 0.1214983946015601
$ws.Range("B5").Value = 0.8580241928393715
$ws.Range("C5").Value = 0.7519232406493396

$ws.Range("A6").Value = 0.09071705834422165
$ws.Range("B6").Value = 0.6381327816132304
$ws.Range("C6").Value = 0.7510679062520823

$ws.Range("A7").Value = 0.1336345653602381
$ws.Range("B7").Value = 0.4060557723420111
$ws.Range("C7").Value = 0.5047731781554882

$ws.Range("A10").Value = 0.1927249413208006
$ws.Range("B10").Value = 0.7091031794084663
$ws.Range("C10").Value = 0.5725905260861619
